$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.073.72"

$ws.Range("D3").Value = "1.962.52"
$ws.Range("E3").Value = "  -6.44%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.95%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.80%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4975"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4207"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09038"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.097"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.34%  "

$ws.Range("D13").Value = "1.959.13"
$ws.Range("E13").Value = "  -7.49%  "

$ws.Range("E14").Value = "  -8.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.423"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.50%  "

$ws.Range("E16").Value = "  +0.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001100"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -10.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06656"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.12%  "

$ws.Range("E20").Value = "  -9.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.920"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.55%  "

$ws.Range("D23").Value = "29.075.90"
$ws.Range("E23").Value = "  -3.97%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.292"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.79%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.58%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.205"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -11.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.255"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.42%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.02%  "

$ws.Range("E31").Value = "  -8.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09838"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.522"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.17%  "

$ws.Range("E34").Value = "  -7.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.706"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.43%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02421"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.992"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06330"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.283"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6422"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.85%  "

$ws.Range("E41").Value = "  -9.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1994"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.83%  "

$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6192"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.168"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.286"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.471"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000330"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06861"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.103"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.86%  "
